$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) values in
# column F for a number of existing rows. No rows inserted on this sheet.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoF = @{2=493; 3=1481; 4=769; 5=185; 7=1072; 8=647; 9=739; 10=1286; 12=996; 13=49; 15=35; 16=383; 18=283; 19=520; 20=537; 23=150; 24=356}
foreach ($row in $expoF.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoF[$row]
}

# ---------------------------------------------------------------------------
# Sheet "演出" (Performance)
# ---------------------------------------------------------------------------
$wsPerf = $wb.Worksheets.Item("演出")

# Update F column counts for existing rows first (rows untouched by the insert).
$perfF = @{3=964; 5=211; 9=560}
foreach ($row in $perfF.Keys) {
    $wsPerf.Cells.Item($row, 6).Value = $perfF[$row]
}

# Insert a new row at position 10, pushing the old row 10 ("Kyle Xian") down to
# row 11 and the old row 11 ("夏川里美") down to row 12.
$wsPerf.Cells.Item(10, 1).EntireRow.Insert()

# Fix up the serial-number column (A) for the two rows that shifted down, so it
# keeps matching (row number - 1) as it does for every other row in the sheet.
$wsPerf.Cells.Item(11, 1).Value = 10
$wsPerf.Cells.Item(12, 1).Value = 11

# Populate the brand-new row 10 with the new event's data.
$wsPerf.Cells.Item(10, 1).Value = 9
$wsPerf.Cells.Item(10, 2).NumberFormat = "@"
$wsPerf.Cells.Item(10, 2).Value = "2024-03-31"
$wsPerf.Cells.Item(10, 3).Value = "【大会员抢先购】广州·KANAKO ITO&AYANE 2024 LIVE"
$wsPerf.Cells.Item(10, 4).Value = "奥体南路12号优托邦购物中心 疆进酒Omni Space GZ"
$wsPerf.Cells.Item(10, 5).Value = "2024.03.31 19:00-03.31 20:30"
$wsPerf.Cells.Item(10, 6).Value = 1
$wsPerf.Cells.Item(10, 7).Value = 380
$wsPerf.Cells.Item(10, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81422"
$wsPerf.Cells.Item(10, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/4Y4U8tC01706172039039.jpeg"

# Restore plain (unstyled) formatting on the new row by copying the format of
# the row directly below it (which still carries the sheet's normal look).
$wsPerf.Range($wsPerf.Cells.Item(11, 1), $wsPerf.Cells.Item(11, 9)).Copy()
$wsPerf.Range($wsPerf.Cells.Item(10, 1), $wsPerf.Cells.Item(10, 9)).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Sheet "全部类型" (All types) - same structural edit as "演出" above, mirrored
# onto this aggregate sheet (rows 34-36 here correspond to rows 9-11 there).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")

# Update F column counts for existing rows first.
$allF = @{2=493; 4=1481; 6=769; 7=185; 8=964; 10=1072; 11=647; 12=739; 13=1286; 15=996; 16=49; 18=35; 19=383; 21=211; 23=283; 27=520; 28=537; 32=150; 33=560}
foreach ($row in $allF.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allF[$row]
}

# Insert a new row at position 34, pushing old row 34 ("Kyle Xian") to row 35,
# old row 35 ("Arknights Only") to row 36, and old row 36 ("夏川里美") to row 37.
$wsAll.Cells.Item(34, 1).EntireRow.Insert()

# Fix up the serial-number column (A) for the three rows that shifted down.
$wsAll.Cells.Item(35, 1).Value = 34
$wsAll.Cells.Item(36, 1).Value = 35
$wsAll.Cells.Item(37, 1).Value = 36

# The "Arknights Only" event (now on row 36) also had its want-to-go count
# refreshed from 354 to 356.
$wsAll.Cells.Item(36, 6).Value = 356

# Populate the brand-new row 34 with the new event's data.
$wsAll.Cells.Item(34, 1).Value = 33
$wsAll.Cells.Item(34, 2).NumberFormat = "@"
$wsAll.Cells.Item(34, 2).Value = "2024-03-31"
$wsAll.Cells.Item(34, 3).Value = "【大会员抢先购】广州·KANAKO ITO&AYANE 2024 LIVE"
$wsAll.Cells.Item(34, 4).Value = "奥体南路12号优托邦购物中心 疆进酒Omni Space GZ"
$wsAll.Cells.Item(34, 5).Value = "2024.03.31 19:00-03.31 20:30"
$wsAll.Cells.Item(34, 6).Value = 1
$wsAll.Cells.Item(34, 7).Value = 380
$wsAll.Cells.Item(34, 8).Value = "https://show.bilibili.com/platform/detail.html?id=81422"
$wsAll.Cells.Item(34, 9).Value = "//i0.hdslb.com/bfs/openplatform/202401/4Y4U8tC01706172039039.jpeg"

# Restore plain (unstyled) formatting on the new row.
$wsAll.Range($wsAll.Cells.Item(35, 1), $wsAll.Cells.Item(35, 9)).Copy()
$wsAll.Range($wsAll.Cells.Item(34, 1), $wsAll.Cells.Item(34, 9)).PasteSpecial(-4122)
